$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Bring down the same look as the existing 27-Apr table (rows 504:517):
# column A = left/top aligned wrapped "Normal 2" text style with a box border,
# columns B:F = the same style but with an integer number format.
# (Using the single, un-mixed source cells A504 / B504:F504 rather than the
# whole old block avoids carrying over a couple of stray cells further down
# the old table that were mis-styled.)
$ws.Range("A504").Copy()
[void]$ws.Range("A519:A532").PasteSpecial(-4122)
$ws.Range("B504:F504").Copy()
[void]$ws.Range("B519:F532").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Rows.Item(519).RowHeight = 14.7
$ws.Cells.Item(519, 1).Value = "Ananthapur"
$ws.Cells.Item(519, 2).Value = 0
$ws.Cells.Item(519, 3).Value = 53
$ws.Cells.Item(519, 4).Value = 35
$ws.Cells.Item(519, 5).Value = 14
$ws.Cells.Item(519, 6).Value = 4

$ws.Rows.Item(520).RowHeight = 14.7
$ws.Cells.Item(520, 1).Value = "Chittoor"
$ws.Cells.Item(520, 2).Value = 0
$ws.Cells.Item(520, 3).Value = 73
$ws.Cells.Item(520, 4).Value = 57
$ws.Cells.Item(520, 5).Value = 16
$ws.Cells.Item(520, 6).Value = 0

$ws.Rows.Item(521).RowHeight = 14.7
$ws.Cells.Item(521, 1).Value = "East Godavari"
$ws.Cells.Item(521, 2).Value = 0
$ws.Cells.Item(521, 3).Value = 39
$ws.Cells.Item(521, 4).Value = 27
$ws.Cells.Item(521, 5).Value = 12
$ws.Cells.Item(521, 6).Value = 0

$ws.Rows.Item(522).RowHeight = 14.7
$ws.Cells.Item(522, 1).Value = "Guntur"
$ws.Cells.Item(522, 2).Value = 23
$ws.Cells.Item(522, 3).Value = 237
$ws.Cells.Item(522, 4).Value = 200
$ws.Cells.Item(522, 5).Value = 29
$ws.Cells.Item(522, 6).Value = 8

$ws.Rows.Item(523).RowHeight = 14.7
$ws.Cells.Item(523, 1).Value = "Kadapa"
$ws.Cells.Item(523, 2).Value = 0
$ws.Cells.Item(523, 3).Value = 58
$ws.Cells.Item(523, 4).Value = 30
$ws.Cells.Item(523, 5).Value = 28
$ws.Cells.Item(523, 6).Value = 0

$ws.Rows.Item(524).RowHeight = 14.7
$ws.Cells.Item(524, 1).Value = "Krishna"
$ws.Cells.Item(524, 2).Value = 33
$ws.Cells.Item(524, 3).Value = 210
$ws.Cells.Item(524, 4).Value = 173
$ws.Cells.Item(524, 5).Value = 29
$ws.Cells.Item(524, 6).Value = 8

$ws.Rows.Item(525).RowHeight = 14.7
$ws.Cells.Item(525, 1).Value = "Kurnool"
$ws.Cells.Item(525, 2).Value = 13
$ws.Cells.Item(525, 3).Value = 292
$ws.Cells.Item(525, 4).Value = 252
$ws.Cells.Item(525, 5).Value = 31
$ws.Cells.Item(525, 6).Value = 9

$ws.Rows.Item(526).RowHeight = 14.7
$ws.Cells.Item(526, 1).Value = "Nellore"
$ws.Cells.Item(526, 2).Value = 7
$ws.Cells.Item(526, 3).Value = 79
$ws.Cells.Item(526, 4).Value = 54
$ws.Cells.Item(526, 5).Value = 23
$ws.Cells.Item(526, 6).Value = 2

$ws.Rows.Item(527).RowHeight = 14.7
$ws.Cells.Item(527, 1).Value = "Prakasam"
$ws.Cells.Item(527, 2).Value = 0
$ws.Cells.Item(527, 3).Value = 56
$ws.Cells.Item(527, 4).Value = 33
$ws.Cells.Item(527, 5).Value = 23
$ws.Cells.Item(527, 6).Value = 0

$ws.Rows.Item(528).RowHeight = 14.7
$ws.Cells.Item(528, 1).Value = "Srikakulam"
$ws.Cells.Item(528, 2).Value = 1
$ws.Cells.Item(528, 3).Value = 4
$ws.Cells.Item(528, 4).Value = 4
$ws.Cells.Item(528, 5).Value = 0
$ws.Cells.Item(528, 6).Value = 0

$ws.Rows.Item(529).RowHeight = 14.7
$ws.Cells.Item(529, 1).Value = "Vishakapatnam"
$ws.Cells.Item(529, 2).Value = 0
$ws.Cells.Item(529, 3).Value = 22
$ws.Cells.Item(529, 4).Value = 3
$ws.Cells.Item(529, 5).Value = 19
$ws.Cells.Item(529, 6).Value = 0

$ws.Rows.Item(530).RowHeight = 14.7
$ws.Cells.Item(530, 1).Value = "Vizianagaram"
$ws.Cells.Item(530, 2).Value = 0
$ws.Cells.Item(530, 3).Value = 0
$ws.Cells.Item(530, 4).Value = 0
$ws.Cells.Item(530, 5).Value = 0
$ws.Cells.Item(530, 6).Value = 0

$ws.Rows.Item(531).RowHeight = 14.7
$ws.Cells.Item(531, 1).Value = "West Godavari"
$ws.Cells.Item(531, 2).Value = 3
$ws.Cells.Item(531, 3).Value = 54
$ws.Cells.Item(531, 4).Value = 43
$ws.Cells.Item(531, 5).Value = 11
$ws.Cells.Item(531, 6).Value = 0

$ws.Rows.Item(532).RowHeight = 14.7
$ws.Cells.Item(532, 1).Value = "Total"
$ws.Cells.Item(532, 2).Value = 80
$ws.Cells.Item(532, 3).Value = 1177
$ws.Cells.Item(532, 4).Value = 911
$ws.Cells.Item(532, 5).Value = 235
$ws.Cells.Item(532, 6).Value = 31

# Update the sheet selection to reflect where the new data was entered
[void]$ws.Range("A519:F532").Select()
